# Natmi L1cam-Erbb2 edit: add a third cluster ("FAPs") to the sending/target
# cluster set, so the interaction table now covers the full 3x3 grid of
# Sending cluster x Target cluster combinations (ECs, FAPs, sCs) for the
# fixed Ligand symbol "L1cam" / Receptor symbol "Erbb2" pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Row 2..10: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E..T = the numeric NATMI metrics for that Sending/Target cluster pair.
$rows = @(
    ,@("ECs",  "L1cam", "Erbb2", "ECs",  3, 1, 19.72083766666667,   59.162513,         0.8016210077351786,   0.8016210077351787,   2, 0.6666666666666666, 2.179771666666667, 6.539315,           0.2349306639444428, 0.2349306639444428, 42.98692318873279,  386.882308698595,  0.1883253555790388,    0.1883253555790388)
    ,@("ECs",  "L1cam", "Erbb2", "FAPs", 3, 1, 19.72083766666667,   59.162513,         0.8016210077351786,   0.8016210077351787,   3, 1,                  4.237801,           12.713403,          0.4567402255103586, 0.4567402255103586, 83.57298558463768,  752.1568702617392, 0.3661325598468064,    0.3661325598468064)
    ,@("ECs",  "L1cam", "Erbb2", "sCs",  3, 1, 19.72083766666667,   59.162513,         0.8016210077351786,   0.8016210077351787,   3, 1,                  2.860789,           8.582367000000001, 0.3083291105451987, 0.3083291105451986, 56.41715546758568,  507.7543992082711, 0.2471630923093334,    0.2471630923093334)
    ,@("FAPs", "L1cam", "Erbb2", "ECs",  1, 0.3333333333333333, 0.099159,             0.297477,          0.004030657259573097,  0.004030657259573097, 2, 0.6666666666666666, 2.179771666666667, 6.539315,           0.2349306639444428, 0.2349306639444428, 0.216143978695,     1.945295808255,    0.0009469249861239959, 0.000946924986123996)
    ,@("FAPs", "L1cam", "Erbb2", "FAPs", 1, 0.3333333333333333, 0.099159,             0.297477,          0.004030657259573097,  0.004030657259573097, 3, 1,                  4.237801,           12.713403,          0.4567402255103586, 0.4567402255103586, 0.420216109359,     3.781944984231,    0.00184096330569238,   0.00184096330569238)
    ,@("FAPs", "L1cam", "Erbb2", "sCs",  1, 0.3333333333333333, 0.099159,             0.297477,          0.004030657259573097,  0.004030657259573097, 3, 1,                  2.860789,           8.582367000000001, 0.3083291105451987, 0.3083291105451986, 0.283672976451,     2.553056788059,    0.001242768967756721,  0.001242768967756721)
    ,@("sCs",  "L1cam", "Erbb2", "ECs",  3, 1, 4.781202,            14.343606,         0.1943483350052483,   0.1943483350052483,   2, 0.6666666666666666, 2.179771666666667, 6.539315,           0.2349306639444428, 0.2349306639444428, 10.42192865221,     93.79735786989001, 0.04565838337927997,   0.04565838337927997)
    ,@("sCs",  "L1cam", "Erbb2", "FAPs", 3, 1, 4.781202,            14.343606,         0.1943483350052483,   0.1943483350052483,   3, 1,                  4.237801,           12.713403,          0.4567402255103586, 0.4567402255103586, 20.261782616802,    182.356043551218,  0.0887667023578598,    0.08876670235785981)
    ,@("sCs",  "L1cam", "Erbb2", "sCs",  3, 1, 4.781202,            14.343606,         0.1943483350052483,   0.1943483350052483,   3, 1,                  2.860789,           8.582367000000001, 0.3083291105451987, 0.3083291105451986, 13.678010088378,    123.102090795402,  0.0599232492681085,    0.05992324926810849)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 2
    $rowVals = $rows[$r]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$rowNum").Value = $rowVals[$i]
    }
}
